$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53, shifting existing rows 53:58 down to 54:59.
$ws.Rows.Item(53).Insert()

# Fill the new row 53 with its data (copy of the constant columns + the new values).
$ws.Cells.Item(53, 1).Value = 2
$ws.Cells.Item(53, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44748
$ws.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 5).Value = 4
$ws.Cells.Item(53, 6).Value = 100112022
$ws.Cells.Item(53, 7).Value = "Arveja Verde"
$ws.Cells.Item(53, 8).Value = "Perfection"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 700
$ws.Cells.Item(53, 11).Value = 28000
$ws.Cells.Item(53, 12).Value = 30000
$ws.Cells.Item(53, 13).Value = 29000
$ws.Cells.Item(53, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(53, 16).Value = 1160
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
